# Update cryptocurrency price (column D) and 1h volume change (column E) values.
# Values that look like plain numbers are entered with a leading apostrophe so
# Excel keeps them as text (matching the original text-typed cell contents)
# instead of auto-converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.945.81'
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').Value = '1.785.97'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''221.73'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '''31.43'
$ws.Range('E8').Value = '  -4.11%  '
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('E10').Value = '  +5.57%  '
$ws.Range('D11').Value = '''0.0922'
$ws.Range('E11').Value = '  -1.61%  '
$ws.Range('D12').Value = '2.042.93'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').Value = '1.787.03'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').Value = '''10.56'
$ws.Range('E14').Value = '  -5.53%  '
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('D16').Value = '33.956.72'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('D18').Value = '''67.97'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('E19').Value = '  -3.11%  '
$ws.Range('D20').Value = '0.0₃0782'
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').Value = '''10.72'
$ws.Range('E22').Value = '  +2.98%  '
$ws.Range('E23').Value = '  -3.33%  '
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('D25').Value = '''157.99'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('D27').Value = '''6.99'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('E28').Value = '  -1.87%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '''0.0518'
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('E32').Value = '  -1.50%  '
$ws.Range('D33').Value = '''3.49'
$ws.Range('E33').Value = '  -1.97%  '
$ws.Range('E34').Value = '  -1.96%  '
$ws.Range('D35').Value = '1.410.85'
$ws.Range('E35').Value = '  -2.06%  '
$ws.Range('E36').Value = '  +1.93%  '
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('E38').Value = '  -1.52%  '
$ws.Range('D39').Value = '''0.937'
$ws.Range('E39').Value = '  +4.16%  '
$ws.Range('D40').Value = '''79.53'
$ws.Range('E40').Value = '  -4.07%  '
$ws.Range('D41').Value = '''2.72'
$ws.Range('E41').Value = '  -3.11%  '
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('D44').Value = '''5.93'
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('D46').Value = '1.941.61'
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('E47').Value = '  -1.02%  '
$ws.Range('D48').Value = '''105.47'
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').Value = '''11.82'
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').Value = '0.0₆0119'
$ws.Range('E51').Value = '  -1.51%  '
